# Updates the cryptos list: refreshed Price (col D) and Volume(1h) (col E)
# values for most rows, plus a reorder/update of the RenderToken /
# LidoDAOToken rows (36 and 37 swapped places with new data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.588.41'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '1.960.73'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.56'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.66'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.377'
$ws.Range("E9").Value = '  +2.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0807'
$ws.Range("E10").Value = '  -6.17%  '
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.21'
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.831'
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").Value = '2.248.82'
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.74'
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("D17").Value = '1.972.85'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").Value = '36.508.56'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.80'
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("D20").Value = '0.0₃0855'
$ws.Range("E20").Value = '  -2.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '228.76'
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.05'
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.45'
$ws.Range("E24").Value = '  -1.40%  '
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.25'
$ws.Range("E26").Value = '  -1.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.138'
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.54'
$ws.Range("E28").Value = '  -1.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.43'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'
$ws.Range("E31").Value = '  -3.39%  '
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.32'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.40'
$ws.Range("E36").Value = '  +11.48%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.24'
$ws.Range("E37").Value = '  +2.15%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.72'
$ws.Range("E39").Value = '  -11.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0981'
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("E42").Value = '  -1.48%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.02'
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("D45").Value = '1.366.96'
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("E46").Value = '  -0.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.95'
$ws.Range("E47").Value = '  -0.63%  '
$ws.Range("E48").Value = '  -1.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.82'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").Value = '2.139.84'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("E51").Value = '  -5.02%  '
